$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Weekly crime statistics updates ---
# Row 14
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 2

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 71.428571428571
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = -25
$ws.Range("N16").Value = -53.846153846153

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -5.882352941176
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 8.333333333333
$ws.Range("L17").Value = 116.666666666667
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = -27.777777777777

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = -50
$ws.Range("L18").Value = 100
$ws.Range("M18").Value = -69.230769230769
$ws.Range("N18").Value = -92.452830188679

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -26.190476190476
$ws.Range("I19").Value = 24
$ws.Range("J19").Value = 34
$ws.Range("K19").Value = -29.411764705882
$ws.Range("L19").Value = -36.842105263157
$ws.Range("M19").Value = -11.111111111111
$ws.Range("N19").Value = -38.461538461538

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 20
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -33.333333333333
$ws.Range("N20").Value = -97.235023041474

# Row 21
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -45.833333333333
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -13.793103448275
$ws.Range("I21").Value = 59
$ws.Range("J21").Value = 68
$ws.Range("K21").Value = -13.235294117647
$ws.Range("L21").Value = -3.278688524590
$ws.Range("M21").Value = -16.901408450704
$ws.Range("N21").Value = -83.473389355742

# Row 23
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 300
$ws.Range("M23").Value = -33.333333333333

# Row 24
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = -9.756097560975
$ws.Range("I24").Value = 54
$ws.Range("J24").Value = 63
$ws.Range("K24").Value = -14.285714285714
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -15.625

# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 50
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 16.666666666666
$ws.Range("I25").Value = 15
$ws.Range("J25").Value = 12
$ws.Range("K25").Value = 25
$ws.Range("L25").Value = 66.666666666666
$ws.Range("M25").Value = -25

# Row 26
$ws.Range("D26").Value = 2
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 2
$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("K26").Value = -100
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = -50

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = -66.666666666666
$ws.Range("M28").Value = 0
$ws.Range("M28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N28").Value = 0
$ws.Range("N28").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = -50
$ws.Range("M29").Value = 0
$ws.Range("M29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N29").Value = 0
$ws.Range("N29").NumberFormat = '#,##0.0;"-"#,##0.0'
